{"js": "// Resume edits:\n//   1. \"Made 2-D Javascript simulation of the International Aerial Robotics\n//      Competition (IARC)\" -> \"Innovated to create 2-D Javascript simulation\n//      of the International Aerial Robotics Competition (IARC)\"\n//   2. \"Participant\" (ShamHacks role) -> \"Competitor\"\n\nconst body = context.document.body;\n\n// 1) Multi-rotor Robot Design Team bullet: \"Made\" -> \"Innovated to create\"\n//    The original text is split by a \"_GoBack\" bookmark:\n//      \"Made 2-D Javascript simulation of the International Aerial Robotics\n//       Competitio\" + <bookmark> + \"n (IARC)\"\n//    Replace on each side of the bookmark separately so the bookmark stays\n//    put instead of being swallowed by a single edge-to-edge replacement.\nconst madeResults = body.search(\n  \"Made 2-D Javascript simulation of the International Aerial Robotics Competitio\",\n  { matchCase: true }\n);\nmadeResults.load(\"items\");\nawait context.sync();\n\nif (madeResults.items.length > 0) {\n  madeResults.items[0].insertText(\"Innovated to create\", \"Replace\");\n  await context.sync();\n} else {\n  // Fall back in case the run boundaries prevent matching the full phrase\n  // (search the shorter, unique lead-in word instead).\n  const fallback = body.search(\"Made 2-D \", { matchCase: true });\n  fallback.load(\"items\");\n  await context.sync();\n  if (fallback.items.length > 0) {\n    fallback.items[0].insertText(\"Innovated to create 2-D \", \"Replace\");\n    await context.sync();\n  }\n}\n\nconst iarcResults = body.search(\"n (IARC)\", { matchCase: true });\niarcResults.load(\"items\");\nawait context.sync();\n\nif (iarcResults.items.length > 0) {\n  iarcResults.items[0].insertText(\n    \" 2-D Javascript simulation of the International Aerial Robotics Competition (IARC)\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) ShamHacks entry: \"Participant\" -> \"Competitor\"\nconst participantResults = body.search(\"Participant\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\nparticipantResults.load(\"items\");\nawait context.sync();\n\nif (participantResults.items.length > 0) {\n  participantResults.items[0].insertText(\"Competitor\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Resume edits:\n#   1. \"Made 2-D Javascript simulation of the International Aerial Robotics\n#      Competition (IARC)\" -> \"Innovated to create 2-D Javascript simulation\n#      of the International Aerial Robotics Competition (IARC)\"\n#   2. \"Participant\" (ShamHacks role) -> \"Competitor\"\n\n$d = $word.ActiveDocument\n\n# 1) Multi-rotor Robot Design Team bullet: \"Made\" -> \"Innovated to create\"\n$oldIarc = \"Made 2-D Javascript simulation of the International Aerial Robotics Competition (IARC)\"\n$newIarc = \"Innovated to create 2-D Javascript simulation of the International Aerial Robotics Competition (IARC)\"\n\n$range1 = $d.Content\n$found1 = $range1.Find.Execute($oldIarc, $false, $false, $false, $false, $false, $true, 1, $false, $newIarc, 2)\n\nif (-not $found1) {\n    # Fallback in case run boundaries keep the whole phrase from matching.\n    $range1b = $d.Content\n    $range1b.Find.Execute(\"Made 2-D \", $false, $false, $false, $false, $false, $true, 1, $false, \"Innovated to create 2-D \", 2)\n}\n\n# 2) ShamHacks entry: \"Participant\" -> \"Competitor\"\n$range2 = $d.Content\n$range2.Find.Execute(\"Participant\", $false, $true, $false, $false, $false, $true, 1, $false, \"Competitor\", 2)\n\n$d.Save()\n"}
